$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the AF column (doctor_MA / average) values for rows 4-13
$ws.Range("AF4").Value = 0.6879999999999999
$ws.Range("AF5").Value = 0.9379999999999999
$ws.Range("AF6").Value = 0.794
$ws.Range("AF7").Value = 0.874
$ws.Range("AF8").Value = 0.864
$ws.Range("AF9").Value = 0.75
$ws.Range("AF10").Value = 0.9379999999999999
$ws.Range("AF11").Value = 0.9379999999999999
$ws.Range("AF12").Value = 1.233
$ws.Range("AF13").Value = 1.688
